$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-08-14"

# Update the August row label shared string
$ws.Range("A9").Value = "August (through 08-14)"

# Update the August (row 9) values for each year column (B..I)
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 31
$ws.Range("E9").Value = 23
$ws.Range("F9").Value = 19
$ws.Range("G9").Value = 86
$ws.Range("H9").Value = 80
$ws.Range("I9").Value = 80

# Update the Total (row 10) values for each year column (B..I)
$ws.Range("B10").Value = 176
$ws.Range("C10").Value = 334
$ws.Range("D10").Value = 496
$ws.Range("E10").Value = 448
$ws.Range("F10").Value = 323
$ws.Range("G10").Value = 707
$ws.Range("H10").Value = 990
$ws.Range("I10").Value = 1050
